$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GDPbES")

# --- New rows 15-17: additional electricity sources that mirror existing
#     guaranteed-dispatch rows (petroleum -> row 11, biomass -> row 9) ---

# Row 15: crude oil (mirrors petroleum, row 11)
$ws.Range("A15").Value = "crude oil"
$ws.Range("B15").Formula = "=B11"
$ws.Range("C15:AK15").Formula = "=C11"

# Row 16: heavy or residual fuel oil (mirrors petroleum, row 11)
$ws.Range("A16").Value = "heavy or residual fuel oil"
$ws.Range("B16").Formula = "=B11"
$ws.Range("C16:AK16").Formula = "=C11"

# Row 17: municipal solid waste (mirrors biomass, row 9)
$ws.Range("A17").Value = "municipal solid waste"
$ws.Range("B17").Formula = "=B9"
$ws.Range("C17:AK17").Formula = "=C9"

# --- Row 1: add a header label in A1 describing the units, bold + wrapped,
#     with the row grown tall enough to show the wrapped text ---
$ws.Range("A1").Value = "Guaranteed Dispatch Fraction (dimensionless)"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 45
